$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 733
$ws.Range("I125").Value = 703.75
$ws.Range("J125").Value = 850
$ws.Range("K125").Value = 6333.75
$ws.Range("L125").Value = 7650
$ws.Range("M125").Value = -3873.75
$ws.Range("N125").Value = -12570

$ws.Range("H132").Value = 4459987.5
$ws.Range("I132").Value = 6462.4443
$ws.Range("J132").Value = 24500850
$ws.Range("K132").Value = 19387.3329
$ws.Range("L132").Value = 73502550
$ws.Range("M132").Value = -16857.3329
$ws.Range("N132").Value = -73507610

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4684.6665
$ws.Range("I45").Value = 4639.5
$ws.Range("J45").Value = 4910.5
$ws.Range("K45").Value = 4639.5
$ws.Range("L45").Value = 4910.5
$ws.Range("M45").Value = -4262.5
$ws.Range("N45").Value = -5664.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 129995
$ws.Range("J42").Value = 129995
$ws.Range("L42").Value = 129995
$ws.Range("N42").Value = -130651

$ws.Range("H134").Value = 6108.353
$ws.Range("I134").Value = 6108.353
$ws.Range("K134").Value = 18325.059
$ws.Range("M134").Value = -15790.059

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2195.6072
$ws.Range("I31").Value = 1478.4584
$ws.Range("J31").Value = 6498.5
$ws.Range("K31").Value = 1478.4584
$ws.Range("L31").Value = 6498.5
$ws.Range("M31").Value = -1183.4584
$ws.Range("N31").Value = -7088.5

$ws.Range("H34").Value = 2195.6072
$ws.Range("I34").Value = 1478.4584
$ws.Range("J34").Value = 6498.5
$ws.Range("K34").Value = 1478.4584
$ws.Range("L34").Value = 6498.5
$ws.Range("M34").Value = -1276.4584
$ws.Range("N34").Value = -6902.5

$ws.Range("H99").Value = 3714.95
$ws.Range("I99").Value = 3714.95
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3714.95
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2216.95
$ws.Range("N99").ClearContents()

$ws.Range("H126").Value = 3714.95
$ws.Range("I126").Value = 3714.95
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11144.85
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8674.849999999999
$ws.Range("N126").ClearContents()

$ws.Range("H127").Value = 33888.89
$ws.Range("J127").Value = 33888.89
$ws.Range("L127").Value = 33888.89
$ws.Range("N127").Value = -43808.89

$ws.Range("H132").Value = 74837.07000000001
$ws.Range("I132").Value = 3247.3635
$ws.Range("J132").Value = 337332.66
$ws.Range("K132").Value = 9742.0905
$ws.Range("L132").Value = 1011997.98
$ws.Range("M132").Value = -7212.0905
$ws.Range("N132").Value = -1017057.98

$ws.Range("H134").Value = 48263.652
$ws.Range("I134").Value = 1284.2667
$ws.Range("J134").Value = 136350
$ws.Range("K134").Value = 3852.800099999999
$ws.Range("L134").Value = 409050
$ws.Range("M134").Value = -1317.800099999999
$ws.Range("N134").Value = -414120

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 7027
$ws.Range("J95").Value = 7027
$ws.Range("L95").Value = 21081
$ws.Range("N95").Value = -25199

$ws.Range("H132").Value = 1508.75
$ws.Range("I132").Value = 1340.125
$ws.Range("J132").Value = 1677.375
$ws.Range("K132").Value = 12061.125
$ws.Range("L132").Value = 15096.375
$ws.Range("M132").Value = -9531.125
$ws.Range("N132").Value = -20156.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13688
$ws.Range("I80").Value = 19916.666
$ws.Range("J80").Value = 4345
$ws.Range("K80").Value = 19916.666
$ws.Range("L80").Value = 4345
$ws.Range("M80").Value = -18918.666
$ws.Range("N80").Value = -6341

$ws.Range("H83").Value = 13688
$ws.Range("I83").Value = 19916.666
$ws.Range("J83").Value = 4345
$ws.Range("K83").Value = 99583.33
$ws.Range("L83").Value = 21725
$ws.Range("M83").Value = -94591.33
$ws.Range("N83").Value = -31709

$ws.Range("H102").Value = 1742.421
$ws.Range("I102").Value = 1306.5454
$ws.Range("J102").Value = 2341.75
$ws.Range("K102").Value = 1306.5454
$ws.Range("L102").Value = 2341.75
$ws.Range("M102").Value = 315.4546
$ws.Range("N102").Value = -5585.75

$ws.Range("H113").Value = 1813.6
$ws.Range("I113").Value = 1638
$ws.Range("K113").Value = 1638
$ws.Range("M113").Value = 532

$ws.Range("H122").Value = 2733.1667
$ws.Range("I122").Value = 1150
$ws.Range("J122").Value = 4316.3335
$ws.Range("K122").Value = 3450
$ws.Range("L122").Value = 12949.0005
$ws.Range("M122").Value = -1000
$ws.Range("N122").Value = -17849.0005

$ws.Range("H126").Value = 1991
$ws.Range("I126").Value = 1800
$ws.Range("J126").Value = 2118.3333
$ws.Range("K126").Value = 5400
$ws.Range("L126").Value = 6354.999899999999
$ws.Range("M126").Value = -2930
$ws.Range("N126").Value = -11294.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2584.3809
$ws.Range("I7").Value = 2531.8333
$ws.Range("J7").Value = 2899.6667
$ws.Range("K7").Value = 2531.8333
$ws.Range("L7").Value = 2899.6667
$ws.Range("M7").Value = -2419.8333
$ws.Range("N7").Value = -3123.6667

$ws.Range("H34").Value = 50012
$ws.Range("J34").Value = 50012
$ws.Range("L34").Value = 50012
$ws.Range("N34").Value = -50356

$ws.Range("H36").Value = 89997.5
$ws.Range("J36").Value = 89997.5
$ws.Range("L36").Value = 89997.5
$ws.Range("N36").Value = -91121.5

$ws.Range("H40").Value = 6714.2856
$ws.Range("I40").Value = 6714.2856
$ws.Range("K40").Value = 6714.2856
$ws.Range("M40").Value = -6578.2856

$ws.Range("H122").Value = 3976
$ws.Range("I122").Value = 3873.3333
$ws.Range("J122").Value = 4900
$ws.Range("K122").Value = 11619.9999
$ws.Range("L122").Value = 14700
$ws.Range("M122").Value = -9169.999899999999
$ws.Range("N122").Value = -19600

$ws.Range("H126").Value = 2584.3809
$ws.Range("I126").Value = 2531.8333
$ws.Range("J126").Value = 2899.6667
$ws.Range("K126").Value = 7595.499899999999
$ws.Range("L126").Value = 8699.000100000001
$ws.Range("M126").Value = -5125.499899999999
$ws.Range("N126").Value = -13639.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2469.4119
$ws.Range("I122").Value = 2165
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 6495
$ws.Range("L122").Value = 9600
$ws.Range("M122").Value = -4045
$ws.Range("N122").Value = -14500

$ws.Range("H126").Value = 1565.8572
$ws.Range("I126").Value = 1076.8
$ws.Range("K126").Value = 3230.4
$ws.Range("M126").Value = -760.3999999999996
